$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cell = $ws.Cells.Item(6,3)
$v = $cell.Value
Write-Output "TYPE=$($cell.GetType())"
Write-Output "VALUE=$v"
Write-Output "VALUE2=$($cell.Value)"
